# Auto-generated edit script implementing the "fixed a few issues with ratios dataframe" commit
# Applies value + column-width changes to the Calc and Results sheets.
$wb = $excel.ActiveWorkbook

# --- Calc sheet ---
$calc = $wb.Worksheets.Item("Calc")

$calc.Range("AT3").Value = 8
$calc.Range("AU3").Value = 4
$calc.Range("AP4").Value = 0.007900000000000001
$calc.Range("AQ4").Value = 1.853155055125499
$calc.Range("AT4").Value = 8
$calc.Range("AU4").Value = 4
$calc.Range("AV4").Value = 0.4213
$calc.Range("AW4").Value = 0.0071802908
$calc.Range("AX4").Value = 0.007376077534584326
$calc.Range("AY4").Value = 1.704317778305246
$calc.Range("BB4").Value = 12.83507455944778
$calc.Range("BC4").Value = 0.6679742161634127
$calc.Range("BD4").Value = 363.3
$calc.Range("BE4").Value = 3.688038767292163
$calc.Range("BF4").Value = 1.750789825441331
$calc.Range("AT5").Value = 8
$calc.Range("AU5").Value = 4
$calc.Range("AP6").Value = 0.008399999999999999
$calc.Range("AQ6").Value = 1.015351142270035
$calc.Range("AT6").Value = 8
$calc.Range("AU6").Value = 4
$calc.Range("AV6").Value = 0.8114
$calc.Range("AW6").Value = 0.0113479214
$calc.Range("AX6").Value = 0.01131192708832987
$calc.Range("AY6").Value = 1.398560685235396
$calc.Range("BB6").Value = 10.61276535539819
$calc.Range("BC6").Value = 0.8522932191132148
$calc.Range("BD6").Value = 753.4
$calc.Range("BE6").Value = 5.655963544164936
$calc.Range("BF6").Value = 1.394124610343834
$calc.Range("AT7").Value = 8
$calc.Range("AU7").Value = 4
$calc.Range("AV7").Value = 683.284
$calc.Range("AX7").Value = 204.7331390773551
$calc.Range("BB7").Value = -3.206491539234178
$calc.Range("BD7").Value = 683226
$calc.Range("BE7").Value = 102366.5695386775
$calc.Range("BF7").Value = 29.96311037245934
$calc.Range("AP8").Value = 0.009299999999999999
$calc.Range("AQ8").Value = 1.047769265434881
$calc.Range("AT8").Value = 8
$calc.Range("AU8").Value = 4
$calc.Range("AV8").Value = 0.8713
$calc.Range("AW8").Value = 0.0122165271
$calc.Range("AX8").Value = 0.01209915221237681
$calc.Range("AY8").Value = 1.402103420176747
$calc.Range("BB8").Value = 10.37093291533475
$calc.Range("BC8").Value = 0.661583215145599
$calc.Range("BD8").Value = 813.3
$calc.Range("BE8").Value = 6.049576106188407
$calc.Range("BF8").Value = 1.388632183217814
$calc.Range("AT9").Value = 8
$calc.Range("AU9").Value = 4
$calc.Range("AT10").Value = 8
$calc.Range("AU10").Value = 4
$calc.Range("AV10").Value = 1.0689
$calc.Range("AW10").Value = 0.009554396600000001
$calc.Range("AX10").Value = 0.009445363760939006
$calc.Range("AY10").Value = 0.89385317616241
$calc.Range("BB10").Value = 12.31164657147083
$calc.Range("BC10").Value = 0.6271336428673672
$calc.Range("BD10").Value = 1010.9
$calc.Range("BE10").Value = 4.722681880469503
$calc.Range("BF10").Value = 0.8836527047374877
$calc.Range("AT11").Value = 8
$calc.Range("AU11").Value = 4
$calc.Range("AP12").Value = 0.0358
$calc.Range("AQ12").Value = 2.129938124702523
$calc.Range("AT12").Value = 8
$calc.Range("AU12").Value = 4
$calc.Range("AV12").Value = 1.1929
$calc.Range("AW12").Value = 0.2488815794
$calc.Range("AX12").Value = 0.247089538595463
$calc.Range("AY12").Value = 20.86357443205633
$calc.Range("BB12").Value = 16.02772876672674
$calc.Range("BC12").Value = 1.715967661765083
$calc.Range("BD12").Value = 1134.9
$calc.Range("BE12").Value = 123.5447692977315
$calc.Range("BF12").Value = 20.71334886373234
$calc.Range("AT13").Value = 8
$calc.Range("AU13").Value = 4
$calc.Range("AP14").Value = 0.0182
$calc.Range("AQ14").Value = 1.202669662327364
$calc.Range("AT14").Value = 8
$calc.Range("AU14").Value = 4
$calc.Range("AV14").Value = 1.4955
$calc.Range("AW14").Value = 0.0199312423
$calc.Range("AX14").Value = 0.02018910639048038
$calc.Range("AY14").Value = 1.332747729856235
$calc.Range("BB14").Value = 16.78624828561402
$calc.Range("BC14").Value = 1.185774376371707
$calc.Range("BD14").Value = 1437.5
$calc.Range("BE14").Value = 10.09455319524019
$calc.Range("BF14").Value = 1.349990397223697
$calc.Range("AT15").Value = 8
$calc.Range("AU15").Value = 4
$calc.Range("AP16").Value = 0.007
$calc.Range("AQ16").Value = 1.071319253137435
$calc.Range("AT16").Value = 8
$calc.Range("AU16").Value = 4
$calc.Range("AV16").Value = 0.5717
$calc.Range("AW16").Value = 0.0418541857
$calc.Range("AX16").Value = 0.04133340390730402
$calc.Range("AY16").Value = 7.321005020115445
$calc.Range("BB16").Value = 9.247056684056117
$calc.Range("BC16").Value = 0.6283778954866831
$calc.Range("BD16").Value = 513.6999999999999
$calc.Range("BE16").Value = 20.66670195365201
$calc.Range("BF16").Value = 7.229911475827186
$calc.Range("AT17").Value = 8
$calc.Range("AU17").Value = 4

# Column width change: col 54 (BB) 19.7109375 -> 20.7109375
$calc.Columns.Item(54).ColumnWidth = 20.7109375

# --- Results sheet ---
$results = $wb.Worksheets.Item("Results")

# Unit label change: "(ng/g)" -> "(μg/g)"
$results.Range("C2").Value = "(μg/g)"

$results.Range("C3").Value = 13.70118282804834
$results.Range("D3").Value = 0.001449630221586743
$results.Range("C4").Value = 1.186316720206925
$results.Range("D4").Value = 0.00004775060295905169
$results.Range("N4").Value = 0.007900000000000001
$results.Range("O4").Value = 0.4213
$results.Range("P4").Value = 0.0071802908
$results.Range("Q4").Value = 12.83507455944778
$results.Range("R4").Value = 0.6679742161634127
$results.Range("C5").Value = 13.69933853789272
$results.Range("D5").Value = 0.001423106824204759
$results.Range("C6").Value = 0.9801591795303184
$results.Range("D6").Value = 0.0000487301575619043
$results.Range("N6").Value = 0.008399999999999999
$results.Range("O6").Value = 0.8114
$results.Range("P6").Value = 0.0113479214
$results.Range("Q6").Value = 10.61276535539819
$results.Range("R6").Value = 0.8522932191132148
$results.Range("C7").Value = 13.69915089929962
$results.Range("D7").Value = 0.001614191900461489
$results.Range("O7").Value = 683.284
$results.Range("Q7").Value = -3.206491539234178
$results.Range("C8").Value = 1.050386684779172
$results.Range("D8").Value = 0.0000595081833313798
$results.Range("N8").Value = 0.009299999999999999
$results.Range("O8").Value = 0.8713
$results.Range("P8").Value = 0.0122165271
$results.Range("Q8").Value = 10.37093291533475
$results.Range("R8").Value = 0.661583215145599
$results.Range("C9").Value = 13.69850581832
$results.Range("D9").Value = 0.001772333129394309
$results.Range("C10").Value = 1.131860529726503
$results.Range("D10").Value = 0.00005667916986329242
$results.Range("O10").Value = 1.0689
$results.Range("P10").Value = 0.009554396600000001
$results.Range("Q10").Value = 12.31164657147083
$results.Range("R10").Value = 0.6271336428673672
$results.Range("C11").Value = 13.69797452749289
$results.Range("D11").Value = 0.002200848263931516
$results.Range("C12").Value = 0.5268169011439808
$results.Range("D12").Value = 0.00003435009050826131
$results.Range("N12").Value = 0.0358
$results.Range("O12").Value = 1.1929
$results.Range("P12").Value = 0.2488815794
$results.Range("Q12").Value = 16.02772876672674
$results.Range("R12").Value = 1.715967661765083
$results.Range("C13").Value = 13.69812723072646
$results.Range("D13").Value = 0.001328377354164738
$results.Range("C14").Value = 0.6324106757473839
$results.Range("D14").Value = 0.00003160028408937892
$results.Range("N14").Value = 0.0182
$results.Range("O14").Value = 1.4955
$results.Range("P14").Value = 0.0199312423
$results.Range("Q14").Value = 16.78624828561402
$results.Range("R14").Value = 1.185774376371707
$results.Range("C15").Value = 13.69798554924028
$results.Range("D15").Value = 0.001491959248085548
$results.Range("C16").Value = 1.771327669508882
$results.Range("D16").Value = 0.00008912990943670092
$results.Range("N16").Value = 0.007
$results.Range("O16").Value = 0.5717
$results.Range("P16").Value = 0.0418541857
$results.Range("Q16").Value = 9.247056684056117
$results.Range("R16").Value = 0.6283778954866831
$results.Range("C17").Value = 13.69798554924028
$results.Range("D17").Value = 0.001491959248085548

# Column width changes
$results.Columns.Item(4).ColumnWidth = 23.7109375   # D col 21.7109375 -> 23.7109375
$results.Columns.Item(17).ColumnWidth = 20.7109375  # Q col 19.7109375 -> 20.7109375

Write-Host "Edit complete"
